# 17.1.2 worksheet: extend the year table from 2021 (col R) through 2022 (col S).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year header in S4, formatted like the existing R4 ("2021") header cell.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("S4").Value = 2022

# New data point in S5, formatted like the existing R5 (72) data cell.
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("S5").Value = 76.099999999999994

$excel.CutCopyMode = $false

# Matches the post-edit selection recorded in the saved workbook.
$ws.Range("P8").Select()
